$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")
$ws.Activate()

# --- Append " order By ... LIMIT 100" clauses to the saved Neo4j queries ---

$b2Suffix = "`n order By ss.study_subject_id ASC LIMIT 100"
$b2 = $ws.Range("B2").Value2
if (-not $b2.EndsWith($b2Suffix)) {
    $ws.Range("B2").Value2 = $b2 + $b2Suffix
}

$b3Suffix = "`n order By samp.sample_id ASC LIMIT 100"
$b3 = $ws.Range("B3").Value2
if (-not $b3.EndsWith($b3Suffix)) {
    $ws.Range("B3").Value2 = $b3 + $b3Suffix
}

$b4 = $ws.Range("B4").Value2
$oldTail = "    order by f.file_name"
$newTail = "     order By f.file_name ASC LIMIT 100"
if ($b4.EndsWith($oldTail)) {
    $b4 = $b4.Substring(0, $b4.Length - $oldTail.Length) + $newTail
    $ws.Range("B4").Value2 = $b4
}

# --- Row heights grow by one wrapped line now that the queries are longer ---
$ws.Rows.Item(2).RowHeight = 331.2
$ws.Rows.Item(3).RowHeight = 360

# --- Selection left on B3 (no more frozen/scrolled topLeftCell) ---
$ws.Range("B3").Select()
